# Divsalar & Simon Compare workbook update
# - Adds a new "DS N=2" column (B) and "DS N=3" column (E) on Sheet1, duplicating
#   the original Div&Sim data to begin the Non-Coherent detector N=2 work.
# - Updates the N=3 Non-Coherent detection values (columns D, M, S) to reflect the
#   extra layer of modulo blocks that improves BER.
# - Adds a matching "DS N=3" series to the first scatter chart.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Header row (row 2): rename B2 "Div&Sim" -> "DS N=2", add E2 "DS N=3" ---
$ws.Range("B2").Value = "DS N=2"
$ws.Range("E2").Value = "DS N=3"

# --- Row 3 (SNR = 7) ---
$ws.Range("B3").Formula = "=5*10^-3"
$ws.Range("E3").Formula = "=0.004"

# --- Row 4 (SNR = 8) ---
$ws.Range("D4").Formula = "=0.000846"
$ws.Range("E4").Formula = "=0.0008"
$ws.Range("M4").Formula = "=0.005532"

# --- Row 5 (SNR = 9) ---
$ws.Range("B5").Formula = "=3*10^-4"
$ws.Range("D5").Value2 = 0.00015
$ws.Range("E5").Formula = "=0.0001"
$ws.Range("M5").Formula = "=0.001758"
$ws.Range("S5").Formula = "=0.00134"

# --- Row 6 (SNR = 10) ---
$ws.Range("B6").Formula = "=5*10^-5"
$ws.Range("E6").Formula = "=0.00001"
$ws.Range("M6").Formula = "=0.000438"

# --- Row 7 (SNR = 11, second table) ---
$ws.Range("M7").Formula = "=0.000062"

# --- Row 8 (SNR = 12, second table) ---
$ws.Range("M8").Formula = "=0.000008"
$ws.Range("S8").Formula = "=0.000001"

# --- New chart series "DS N=3" on the first chart (Non-Coherent Detection M=2, N=2,3) ---
$chartObjs = $ws.ChartObjects()
$co = $chartObjs.Item(1)
$chart = $co.Chart
$sc = $chart.SeriesCollection()
$newSeries = $sc.NewSeries()
$newSeries.Name = "DS N=3"
$newSeries.Formula = '=SERIES("DS N=3",Sheet1!$A$3:$A$6,Sheet1!$E$3:$E$6,4)'

# --- Restore the last active selection ---
$ws.Range("F42").Select() | Out-Null
